$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.054745150301199
$ws.Cells.Item(2, 4).Value = 1.060038307084648
$ws.Cells.Item(2, 5).Value = 1.051017841050022
$ws.Cells.Item(2, 6).Value = 1.069763675009067
$ws.Cells.Item(2, 9).Value = 1.046848962490686
$ws.Cells.Item(2, 10).Value = 1.059755638480245
$ws.Cells.Item(2, 11).Value = 1.062766175295604
$ws.Cells.Item(2, 12).Value = 1.053770508784097
$ws.Cells.Item(2, 13).Value = 1.07246531259788
$ws.Cells.Item(2, 14).Value = 1.061260613042558

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.056480951296423
$ws.Cells.Item(3, 4).Value = 1.061415257188752
$ws.Cells.Item(3, 5).Value = 1.05254222078431
$ws.Cells.Item(3, 6).Value = 1.071261835915073
$ws.Cells.Item(3, 9).Value = 1.047362801924152
$ws.Cells.Item(3, 10).Value = 1.061139039670446
$ws.Cells.Item(3, 11).Value = 1.063956051346688
$ws.Cells.Item(3, 12).Value = 1.055105619597533
$ws.Cells.Item(3, 13).Value = 1.073778017217178
$ws.Cells.Item(3, 14).Value = 1.062645978821127

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.057601026073847
$ws.Cells.Item(4, 4).Value = 1.062303520508216
$ws.Cells.Item(4, 5).Value = 1.05352537210779
$ws.Cells.Item(4, 6).Value = 1.072228848981741
$ws.Cells.Item(4, 9).Value = 1.047692452129098
$ws.Cells.Item(4, 10).Value = 1.062030787836578
$ws.Cells.Item(4, 11).Value = 1.064722736993398
$ws.Cells.Item(4, 12).Value = 1.055965810146917
$ws.Cells.Item(4, 13).Value = 1.074624519816927
$ws.Cells.Item(4, 14).Value = 1.06353899337194

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.058071177802028
$ws.Cells.Item(5, 4).Value = 1.062676307395924
$ws.Cells.Item(5, 5).Value = 1.053937930607902
$ws.Cells.Item(5, 6).Value = 1.072634818871317
$ws.Cells.Item(5, 9).Value = 1.047830362737364
$ws.Cells.Item(5, 10).Value = 1.06240487612436
$ws.Cells.Item(5, 11).Value = 1.065044285501689
$ws.Cells.Item(5, 12).Value = 1.056326557430874
$ws.Cells.Item(5, 13).Value = 1.07497970477937
$ws.Cells.Item(5, 14).Value = 1.063913612907999

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.058150076121273
$ws.Cells.Item(6, 4).Value = 1.062738862802626
$ws.Cells.Item(6, 5).Value = 1.054007156919924
$ws.Cells.Item(6, 6).Value = 1.072702950354473
$ws.Cells.Item(6, 9).Value = 1.047853479140419
$ws.Cells.Item(6, 10).Value = 1.062467640487738
$ws.Cells.Item(6, 11).Value = 1.065098230285188
$ws.Cells.Item(6, 12).Value = 1.056387077435534
$ws.Cells.Item(6, 13).Value = 1.075039302054882
$ws.Cells.Item(6, 14).Value = 1.063976466403972

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.057607311102417
$ws.Cells.Item(7, 4).Value = 1.062308504200508
$ws.Cells.Item(7, 5).Value = 1.053530887695718
$ws.Cells.Item(7, 6).Value = 1.072234275765326
$ws.Cells.Item(7, 9).Value = 1.047694297539984
$ws.Cells.Item(7, 10).Value = 1.062035789559289
$ws.Cells.Item(7, 11).Value = 1.064727036536394
$ws.Cells.Item(7, 12).Value = 1.055970633897765
$ws.Cells.Item(7, 13).Value = 1.074629268491525
$ws.Cells.Item(7, 14).Value = 1.063544002197671

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.055332424573299
$ws.Cells.Item(8, 4).Value = 1.060504222774404
$ws.Cells.Item(8, 5).Value = 1.051533688143244
$ws.Cells.Item(8, 6).Value = 1.070270488812621
$ws.Cells.Item(8, 9).Value = 1.047023207730294
$ws.Cells.Item(8, 10).Value = 1.060223877956374
$ws.Cells.Item(8, 11).Value = 1.063168977411418
$ws.Cells.Item(8, 12).Value = 1.054222492688323
$ws.Cells.Item(8, 13).Value = 1.072909555283214
$ws.Cells.Item(8, 14).Value = 1.061729517472451

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.051299307037036
$ws.Cells.Item(9, 4).Value = 1.057303551201375
$ws.Cells.Item(9, 5).Value = 1.047989060631399
$ws.Cells.Item(9, 6).Value = 1.066791166769652
$ws.Cells.Item(9, 9).Value = 1.045818691344029
$ws.Cells.Item(9, 10).Value = 1.057004416891532
$ws.Cells.Item(9, 11).Value = 1.060398165033794
$ws.Cells.Item(9, 12).Value = 1.051113018234491
$ws.Cells.Item(9, 13).Value = 1.069856466467243
$ws.Cells.Item(9, 14).Value = 1.058505484403619

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.04859312222012
$ws.Cells.Item(10, 4).Value = 1.055154747669073
$ws.Cells.Item(10, 5).Value = 1.04560810035456
$ws.Cells.Item(10, 6).Value = 1.064458196013568
$ws.Cells.Item(10, 9).Value = 1.04500058079728
$ws.Cells.Item(10, 10).Value = 1.054839404854404
$ws.Cells.Item(10, 11).Value = 1.058533276157747
$ws.Cells.Item(10, 12).Value = 1.04901971793479
$ws.Cells.Item(10, 13).Value = 1.067805105907259
$ws.Cells.Item(10, 14).Value = 1.056337397801068

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.047416957839678
$ws.Cells.Item(11, 4).Value = 1.054220571063186
$ws.Cells.Item(11, 5).Value = 1.044572685868584
$ws.Cells.Item(11, 6).Value = 1.063444646496548
$ws.Cells.Item(11, 9).Value = 1.044642674707984
$ws.Cells.Item(11, 10).Value = 1.053897316680274
$ws.Cells.Item(11, 11).Value = 1.057721416669485
$ws.Cells.Item(11, 12).Value = 1.048108297299688
$ws.Cells.Item(11, 13).Value = 1.06691290355878
$ws.Cells.Item(11, 14).Value = 1.055393971753672

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.04697940264444
$ws.Cells.Item(12, 4).Value = 1.053873002100352
$ws.Cells.Item(12, 5).Value = 1.044187402620274
$ws.Cells.Item(12, 6).Value = 1.063067650387324
$ws.Cells.Item(12, 9).Value = 1.044509176312132
$ws.Cells.Item(12, 10).Value = 1.053546673121241
$ws.Cells.Item(12, 11).Value = 1.05741918942305
$ws.Cells.Item(12, 12).Value = 1.047768987159557
$ws.Cells.Item(12, 13).Value = 1.066580893171269
$ws.Cells.Item(12, 14).Value = 1.055042830240565

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.047073290589653
$ws.Cells.Item(13, 4).Value = 1.053947583027389
$ws.Cells.Item(13, 5).Value = 1.044270078416862
$ws.Cells.Item(13, 6).Value = 1.063148541034485
$ws.Cells.Item(13, 9).Value = 1.044537837441937
$ws.Cells.Item(13, 10).Value = 1.053621919736206
$ws.Cells.Item(13, 11).Value = 1.057484048561971
$ws.Cells.Item(13, 12).Value = 1.047841805339026
$ws.Cells.Item(13, 13).Value = 1.066652138196813
$ws.Cells.Item(13, 14).Value = 1.055118183714355

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.047380803279252
$ws.Cells.Item(14, 4).Value = 1.054191852697261
$ws.Cells.Item(14, 5).Value = 1.04454085228999
$ws.Cells.Item(14, 6).Value = 1.063413494544569
$ws.Cells.Item(14, 9).Value = 1.044631651071553
$ws.Cells.Item(14, 10).Value = 1.053868346949206
$ws.Cells.Item(14, 11).Value = 1.057696448142295
$ws.Cells.Item(14, 12).Value = 1.048080265574027
$ws.Cells.Item(14, 13).Value = 1.066885471930881
$ws.Cells.Item(14, 14).Value = 1.055364960882263

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.047570181928714
$ws.Cells.Item(15, 4).Value = 1.054342278805641
$ws.Cells.Item(15, 5).Value = 1.044707593843001
$ws.Cells.Item(15, 6).Value = 1.063576671986793
$ws.Cells.Item(15, 9).Value = 1.044689378857138
$ws.Cells.Item(15, 10).Value = 1.054020084333661
$ws.Cells.Item(15, 11).Value = 1.057827225829482
$ws.Cells.Item(15, 12).Value = 1.048227086567202
$ws.Cells.Item(15, 13).Value = 1.067029155749373
$ws.Cells.Item(15, 14).Value = 1.055516913751208

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.048671086715571
$ws.Cells.Item(16, 4).Value = 1.055216666129891
$ws.Cells.Item(16, 5).Value = 1.045676722374558
$ws.Cells.Item(16, 6).Value = 1.06452539000239
$ws.Cells.Item(16, 9).Value = 1.045024256210511
$ws.Cells.Item(16, 10).Value = 1.054901829417153
$ws.Cells.Item(16, 11).Value = 1.058587063844789
$ws.Cells.Item(16, 12).Value = 1.049080099087244
$ws.Cells.Item(16, 13).Value = 1.067864234090614
$ws.Cells.Item(16, 14).Value = 1.056399911013857

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.049360472247215
$ws.Cells.Item(17, 4).Value = 1.055764137809398
$ws.Cells.Item(17, 5).Value = 1.046283429858841
$ws.Cells.Item(17, 6).Value = 1.065119586788771
$ws.Cells.Item(17, 9).Value = 1.045233331833837
$ws.Cells.Item(17, 10).Value = 1.05545367610174
$ws.Cells.Item(17, 11).Value = 1.059062516399995
$ws.Cells.Item(17, 12).Value = 1.049613819640212
$ws.Cells.Item(17, 13).Value = 1.068386989879309
$ws.Cells.Item(17, 14).Value = 1.056952541384033

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.0497621589662
$ws.Cells.Item(18, 4).Value = 1.05608310977367
$ws.Cells.Item(18, 5).Value = 1.046636884524872
$ws.Cells.Item(18, 6).Value = 1.065465848713343
$ws.Cells.Item(18, 9).Value = 1.045354929523637
$ws.Cells.Item(18, 10).Value = 1.05577511386888
$ws.Cells.Item(18, 11).Value = 1.05933942097028
$ws.Cells.Item(18, 12).Value = 1.049924647690451
$ws.Cells.Item(18, 13).Value = 1.06869152430657
$ws.Cells.Item(18, 14).Value = 1.057274435629672

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.04989905305727
$ws.Cells.Item(19, 4).Value = 1.056191810489486
$ws.Cells.Item(19, 5).Value = 1.046757331415234
$ws.Cells.Item(19, 6).Value = 1.065583860721234
$ws.Cells.Item(19, 9).Value = 1.045396331582719
$ws.Cells.Item(19, 10).Value = 1.055884640813075
$ws.Cells.Item(19, 11).Value = 1.059433767614668
$ws.Cells.Item(19, 12).Value = 1.050030550715822
$ws.Cells.Item(19, 13).Value = 1.068795298591552
$ws.Cells.Item(19, 14).Value = 1.05738411811469

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.049286551304526
$ws.Cells.Item(20, 4).Value = 1.055705436511825
$ws.Cells.Item(20, 5).Value = 1.046218380231084
$ws.Cells.Item(20, 6).Value = 1.065055868609881
$ws.Cells.Item(20, 9).Value = 1.045210936514028
$ws.Cells.Item(20, 10).Value = 1.05539451432172
$ws.Cells.Item(20, 11).Value = 1.05901154825833
$ws.Cells.Item(20, 12).Value = 1.049556606456098
$ws.Cells.Item(20, 13).Value = 1.068330942550375
$ws.Cells.Item(20, 14).Value = 1.056893295587499

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.047290267254551
$ws.Cells.Item(21, 4).Value = 1.054119937328798
$ws.Cells.Item(21, 5).Value = 1.044461135155733
$ws.Cells.Item(21, 6).Value = 1.063335486773099
$ws.Cells.Item(21, 9).Value = 1.044604040700156
$ws.Cells.Item(21, 10).Value = 1.053795800002934
$ws.Cells.Item(21, 11).Value = 1.057633920262432
$ws.Cells.Item(21, 12).Value = 1.048010066302999
$ws.Cells.Item(21, 13).Value = 1.066816777817782
$ws.Cells.Item(21, 14).Value = 1.055292310911006

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.046031207122853
$ws.Cells.Item(22, 4).Value = 1.053119742138026
$ws.Cells.Item(22, 5).Value = 1.043352318412592
$ws.Cells.Item(22, 6).Value = 1.062250806278482
$ws.Cells.Item(22, 9).Value = 1.044219240318257
$ws.Cells.Item(22, 10).Value = 1.052786508876925
$ws.Cells.Item(22, 11).Value = 1.056763888285038
$ws.Cells.Item(22, 12).Value = 1.047033244246559
$ws.Cells.Item(22, 13).Value = 1.065861245403453
$ws.Cells.Item(22, 14).Value = 1.05428158647583

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.046699036438805
$ws.Cells.Item(23, 4).Value = 1.053650284697944
$ws.Cells.Item(23, 5).Value = 1.043940504800553
$ws.Cells.Item(23, 6).Value = 1.062826105862064
$ws.Cells.Item(23, 9).Value = 1.044423537756767
$ws.Cells.Item(23, 10).Value = 1.05332194853861
$ws.Cells.Item(23, 11).Value = 1.05722547917354
$ws.Cells.Item(23, 12).Value = 1.047551503237895
$ws.Cells.Item(23, 13).Value = 1.066368129068753
$ws.Cells.Item(23, 14).Value = 1.054817786523251

$ws.Cells.Item(24, 2).Value = 1.019999999999999
$ws.Cells.Item(24, 3).Value = 1.049319954294705
$ws.Cells.Item(24, 4).Value = 1.055731962210596
$ws.Cells.Item(24, 5).Value = 1.046247774678671
$ws.Cells.Item(24, 6).Value = 1.065084661110162
$ws.Cells.Item(24, 9).Value = 1.045221057084446
$ws.Cells.Item(24, 10).Value = 1.055421248358523
$ws.Cells.Item(24, 11).Value = 1.059034579861749
$ws.Cells.Item(24, 12).Value = 1.049582460120746
$ws.Cells.Item(24, 13).Value = 1.068356269100392
$ws.Cells.Item(24, 14).Value = 1.056920067589702

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.052344965759531
$ws.Cells.Item(25, 4).Value = 1.058133597266422
$ws.Cells.Item(25, 5).Value = 1.048908518796414
$ws.Cells.Item(25, 6).Value = 1.067692966741841
$ws.Cells.Item(25, 9).Value = 1.046132724821781
$ws.Cells.Item(25, 10).Value = 1.057839962927484
$ws.Cells.Item(25, 11).Value = 1.061117553813296
$ws.Cells.Item(25, 12).Value = 1.051920413678576
$ws.Cells.Item(25, 13).Value = 1.070648527046284
$ws.Cells.Item(25, 14).Value = 1.059342217010781
